# Add 2022-Q4 data
# 1. Insert a new worksheet "2022-Q4" right after "总计" (the summary sheet),
#    shifting 2022-Q3 .. 2021-Q2 one position later.
# 2. Populate the new sheet with the Q4 fund-holdings detail table, copying
#    the header/row-index cell formatting from the "2022-Q3" sheet so the
#    new sheet matches the look of the others.
# 3. Insert a new row at the top of the "总计" (summary) sheet's data for
#    2022-Q4, pushing the existing quarters down by one row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Step 1: create the new "2022-Q4" sheet right after "总计"
# ---------------------------------------------------------------------------
$summarySheet = $wb.Worksheets.Item(1)
$q4Sheet = $wb.Worksheets.Add($null, $summarySheet)
$q4Sheet.Name = "2022-Q4"

# NOTE: worksheet variables in this host resolve by *position*, not object
# identity, so "2022-Q3" must be re-fetched now that the new sheet pushed it
# from index 2 to index 3.
$q3Sheet = $wb.Worksheets.Item(3)

# Copy the header-row look (bold + thin border + centered) and the
# row-index column look from the existing "2022-Q3" sheet onto the new one.
$q3Sheet.Range("B1:H1").Copy()
$q4Sheet.Range("B1:H1").PasteSpecial(-4122)

$q3Sheet.Range("A2").Copy()
$q4Sheet.Range("A2:A18").PasteSpecial(-4122)

# Header row values
$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($col = 2; $col -le 8; $col++) {
    $q4Sheet.Cells.Item(1, $col).Value = $headers[$col - 2]
}

# Detail rows: index, code, name, scale, stock position, position ratio,
# market value (billion CNY), position rank
$rows = @(
    @(0,  "501054", "东方红睿泽三年定期开放灵活配置混合A",          "106.71", "96.60", "4.12", "4.3965", 5),
    @(1,  "009576", "东方红智远三年持有期混合",                     "66.48",  "93.65", "4.04", "2.6858", 5),
    @(2,  "169104", "东方红睿满沪港深灵活配置混合（LOF）",           "43.50",  "93.71", "5.26", "2.2881", 2),
    @(3,  "010377", "广发价值核心混合A",                            "24.90",  "89.01", "3.98", "0.9910", 8),
    @(4,  "009863", "富国创新趋势股票",                             "30.80",  "92.72", "2.72", "0.8378", 7),
    @(5,  "011153", "华宝新兴消费混合A",                            "6.81",   "91.82", "4.53", "0.3085", 7),
    @(6,  "010378", "广发价值核心混合C",                            "4.47",   "89.01", "3.98", "0.1779", 8),
    @(7,  "014307", "嘉实多元动力混合A",                            "1.56",   "93.22", "2.98", "0.0465", 10),
    @(8,  "010783", "德邦沪港深龙头混合A",                          "0.70",   "82.71", "2.63", "0.0184", 8),
    @(9,  "005493", "鑫元价值精选灵活配置混合A",                    "0.57",   "86.30", "2.82", "0.0161", 8),
    @(10, "011032", "东方红睿泽三年定期开放灵活配置混合C",          "0.34",   "96.60", "4.12", "0.0140", 5),
    @(11, "013897", "德邦港股通成长精选混合型证券投资基金A",        "0.46",   "83.95", "2.67", "0.0123", 8),
    @(12, "013898", "德邦港股通成长精选混合型证券投资基金C",        "0.42",   "83.95", "2.67", "0.0112", 8),
    @(13, "011154", "华宝新兴消费混合C",                            "0.22",   "91.82", "4.53", "0.0100", 7),
    @(14, "014308", "嘉实多元动力混合C",                            "0.33",   "93.22", "2.98", "0.0098", 10),
    @(15, "010784", "德邦沪港深龙头混合C",                          "0.35",   "82.71", "2.63", "0.0092", 8),
    @(16, "005494", "鑫元价值精选灵活配置混合C",                    "0.00",   "86.30", "2.82", 0,        8)
)

$r = 2
foreach ($row in $rows) {
    $q4Sheet.Cells.Item($r, 1).Value = $row[0]

    # Columns B, D, E, F, G hold numeric-looking *text* (fund codes must
    # keep leading zeros, and the scale/position/value figures are stored
    # as plain text in the source data) - force text format so the COM
    # layer doesn't silently coerce them to numbers.
    $bCell = $q4Sheet.Cells.Item($r, 2)
    $bCell.NumberFormat = "@"
    $bCell.Value = $row[1]

    $q4Sheet.Cells.Item($r, 3).Value = $row[2]

    $dCell = $q4Sheet.Cells.Item($r, 4)
    $dCell.NumberFormat = "@"
    $dCell.Value = $row[3]

    $eCell = $q4Sheet.Cells.Item($r, 5)
    $eCell.NumberFormat = "@"
    $eCell.Value = $row[4]

    $fCell = $q4Sheet.Cells.Item($r, 6)
    $fCell.NumberFormat = "@"
    $fCell.Value = $row[5]

    $gValue = $row[6]
    $gCell = $q4Sheet.Cells.Item($r, 7)
    if ($gValue -is [string]) {
        $gCell.NumberFormat = "@"
    }
    $gCell.Value = $gValue

    $q4Sheet.Cells.Item($r, 8).Value = $row[7]

    $r = $r + 1
}

# ---------------------------------------------------------------------------
# Step 2: insert the 2022-Q4 row into the "总计" summary sheet
# ---------------------------------------------------------------------------
$summarySheet.Rows.Item(2).Insert()

# Row 3 (the old row 2, "2022-Q3") still carries the original bold+border
# formatting for column A - copy it onto the newly-inserted blank row 2.
$summarySheet.Range("A3").Copy()
$summarySheet.Range("A2").PasteSpecial(-4122)

$summarySheet.Cells.Item(2, 1).Value = 0
$summarySheet.Cells.Item(2, 2).Value = "2022-Q4"
$summarySheet.Cells.Item(2, 3).Value = 17
$summarySheet.Cells.Item(2, 4).Value = 11.83

# Renumber the A column (0-based index) for the rows that shifted down.
for ($row = 3; $row -le 8; $row++) {
    $summarySheet.Cells.Item($row, 1).Value = $row - 2
}
